$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 'Śląsk Wrocław'
$ws.Range("C8").Value = 'Remis'
$ws.Range("C12").Value = 'Piast Gliwice'
$ws.Range("C14").Value = 'Lechia Gdańsk'
$ws.Range("C17").Value = 'Radomiak Radom'
$ws.Range("C18").Value = 'Remis'
$ws.Range("C19").Value = 'Pogoń Szczecin'
$ws.Range("C21").Value = 'Remis'
$ws.Range("C22").Value = 'Remis'
$ws.Range("C26").Value = 'Radomiak Radom'
$ws.Range("C27").Value = 'Remis'
$ws.Range("C30").Value = 'Remis'
$ws.Range("C34").Value = 'Stal Mielec'
$ws.Range("C35").Value = 'Warta Poznań'
$ws.Range("C36").Value = 'Wisła Płock'
$ws.Range("C37").Value = 'Śląsk Wrocław'
$ws.Range("C39").Value = 'Remis'
$ws.Range("C40").Value = 'Lech Poznań'
$ws.Range("C41").Value = 'Zagłębie Lubin'
$ws.Range("C42").Value = 'Górnik Zabrze'
$ws.Range("C45").Value = 'Jagielonia Białystok'
$ws.Range("C47").Value = 'Piast Gliwice'
$ws.Range("C54").Value = 'Wisła Płock'
$ws.Range("C55").Value = 'Śląsk Wrocław'
$ws.Range("C58").Value = 'Remis'
$ws.Range("C59").Value = 'Remis'
$ws.Range("C63").Value = 'Widzew Łódź'
$ws.Range("C65").Value = 'Remis'
$ws.Range("C66").Value = 'Remis'
$ws.Range("C72").Value = 'Remis'
$ws.Range("C76").Value = 'Remis'
$ws.Range("C78").Value = 'Wisła Płock'
$ws.Range("C79").Value = 'Remis'
$ws.Range("C80").Value = 'Zagłębie Lubin'
$ws.Range("C81").Value = 'Remis'
$ws.Range("C82").Value = 'Śląsk Wrocław'
$ws.Range("C83").Value = 'Pogoń Szczecin'
$ws.Range("C84").Value = 'Remis'
$ws.Range("C90").Value = 'Stal Mielec'
$ws.Range("C95").Value = 'Górnik Zabrze'
$ws.Range("C97").Value = 'Remis'
$ws.Range("C98").Value = 'Raków Częstochowa'
$ws.Range("C99").Value = 'Piast Gliwice'
$ws.Range("C104").Value = 'Zagłębie Lubin'
$ws.Range("C111").Value = 'Lechia Gdańsk'
$ws.Range("C112").Value = 'Remis'
$ws.Range("C113").Value = 'Lech Poznań'
$ws.Range("C115").Value = 'Remis'
$ws.Range("C116").Value = 'Warta Poznań'
$ws.Range("C117").Value = 'Remis'
$ws.Range("C118").Value = 'Legia Warszawa'
$ws.Range("C120").Value = 'Remis'
$ws.Range("C121").Value = 'Remis'
$ws.Range("C123").Value = 'Raków Częstochowa'
$ws.Range("C124").Value = 'Remis'
$ws.Range("C126").Value = 'Widzew Łódź'
$ws.Range("C127").Value = 'Remis'
$ws.Range("C129").Value = 'Zagłębie Lubin'
$ws.Range("C131").Value = 'Remis'
$ws.Range("C132").Value = 'Stal Mielec'
$ws.Range("C134").Value = 'Górnik Zabrze'
$ws.Range("C135").Value = 'Remis'
$ws.Range("C136").Value = 'Remis'
$ws.Range("C137").Value = 'Remis'
$ws.Range("C138").Value = 'Remis'
$ws.Range("C140").Value = 'Lechia Gdańsk'
$ws.Range("C141").Value = 'Miedź Legnica'
$ws.Range("C145").Value = 'Widzew Łódź'
$ws.Range("C147").Value = 'Remis'
$ws.Range("C148").Value = 'Widzew Łódź'
$ws.Range("C149").Value = 'Remis'
$ws.Range("C150").Value = 'Miedź Legnica'
$ws.Range("C156").Value = 'Piast Gliwice'
$ws.Range("C159").Value = 'Remis'
$ws.Range("C160").Value = 'Lech Poznań'
$ws.Range("C163").Value = 'Zagłębie Lubin'
$ws.Range("C164").Value = 'Zagłębie Lubin'
$ws.Range("C165").Value = 'Korona Kielce'
$ws.Range("C167").Value = 'Górnik Zabrze'
$ws.Range("C170").Value = 'Remis'
$ws.Range("C172").Value = 'Warta Poznań'
$ws.Range("C173").Value = 'Pogoń Szczecin'
$ws.Range("C174").Value = 'Piast Gliwice'
$ws.Range("C176").Value = 'Remis'
$ws.Range("C177").Value = 'Remis'
$ws.Range("C180").Value = 'Wisła Płock'
$ws.Range("C181").Value = 'Remis'
$ws.Range("C182").Value = 'Cracovia'
$ws.Range("C183").Value = 'Remis'
$ws.Range("C186").Value = 'Wisła Płock'
$ws.Range("C188").Value = 'Remis'
$ws.Range("C190").Value = 'Śląsk Wrocław'
$ws.Range("C192").Value = 'Raków Częstochowa'
$ws.Range("C193").Value = 'Zagłębie Lubin'
$ws.Range("C194").Value = 'Remis'
$ws.Range("C197").Value = 'Remis'
$ws.Range("C201").Value = 'Wisła Płock'
$ws.Range("C203").Value = 'Miedź Legnica'
$ws.Range("C204").Value = 'Remis'
$ws.Range("C207").Value = 'Stal Mielec'
$ws.Range("C209").Value = 'Remis'
$ws.Range("C210").Value = 'Górnik Zabrze'
$ws.Range("C212").Value = 'Korona Kielce'
$ws.Range("C216").Value = 'Warta Poznań'
$ws.Range("C219").Value = 'Remis'
$ws.Range("C221").Value = 'Remis'
$ws.Range("C222").Value = 'Legia Warszawa'
$ws.Range("C223").Value = 'Raków Częstochowa'
$ws.Range("C226").Value = 'Śląsk Wrocław'
$ws.Range("C228").Value = 'Piast Gliwice'
$ws.Range("C229").Value = 'Zagłębie Lubin'
$ws.Range("C230").Value = 'Korona Kielce'
$ws.Range("C233").Value = 'Legia Warszawa'
$ws.Range("C234").Value = 'Remis'
$ws.Range("C235").Value = 'Wisła Płock'
$ws.Range("C239").Value = 'Górnik Zabrze'
$ws.Range("C242").Value = 'Stal Mielec'
$ws.Range("C246").Value = 'Piast Gliwice'
$ws.Range("C247").Value = 'Zagłębie Lubin'
$ws.Range("C248").Value = 'Remis'
$ws.Range("C250").Value = 'Legia Warszawa'
$ws.Range("C252").Value = 'Miedź Legnica'
$ws.Range("C253").Value = 'Remis'
$ws.Range("C254").Value = 'Zagłębie Lubin'
$ws.Range("C255").Value = 'Cracovia'
$ws.Range("C256").Value = 'Remis'
$ws.Range("C258").Value = 'Stal Mielec'
$ws.Range("C262").Value = 'Wisła Płock'
$ws.Range("C264").Value = 'Piast Gliwice'
$ws.Range("C265").Value = 'Warta Poznań'
$ws.Range("C270").Value = 'Stal Mielec'
$ws.Range("C271").Value = 'Radomiak Radom'
$ws.Range("C272").Value = 'Remis'
$ws.Range("C277").Value = 'Warta Poznań'
$ws.Range("C279").Value = 'Remis'
$ws.Range("C281").Value = 'Remis'
$ws.Range("C288").Value = 'Widzew Łódź'
$ws.Range("C293").Value = 'Remis'
$ws.Range("C294").Value = 'Remis'
$ws.Range("C295").Value = 'Radomiak Radom'
$ws.Range("C296").Value = 'Remis'
$ws.Range("C298").Value = 'Miedź Legnica'
$ws.Range("C299").Value = 'Wisła Płock'
$ws.Range("C302").Value = 'Remis'
$ws.Range("C303").Value = 'Remis'
$ws.Range("C305").Value = 'Raków Częstochowa'
$ws.Range("C306").Value = 'Remis'
